# "add idea of game" -- insert the new "Выбор идеи" section body text.
$d = $word.ActiveDocument

# Locate the "Выбор идеи" heading paragraph (the only paragraph using the custom
# "Мой заголовок 3" style) instead of relying on hard-coded offsets.
$heading = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Style.NameLocal -eq "Мой заголовок 3") {
        $heading = $para
    }
}

# 1) Mark the heading's leading tab run with a lastRenderedPageBreak, matching the
#    page break that now falls right before this heading.
$headingInsertPoint = $d.Range($heading.Range.Start, $heading.Range.Start)
$headingInsertPoint.InsertXML('<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# 2) The heading is followed by a single-space paragraph and then an empty paragraph
#    styled "Мой заголовок 2". Replace both of them with the three new paragraphs
#    that discuss choosing the game's idea (incl. the Fortnite vs. PUBG comparison).
$firstOld = $heading.Next()
$secondOld = $firstOld.Next()
$replaceRange = $d.Range($firstOld.Range.Start, $secondOld.Range.End)
$replaceRange.InsertXML('<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t xml:space="preserve">Прежде, чем будет выбрана основная идея игры, геймдизайнер пройдет не мало игр, </w:t></w:r><w:r><w:t xml:space="preserve">изучит потребительский рынок, прочитает множество форумов и неоднократно поделится своими наработками, как с командой, так и с другими геймдизайнерами. При выборе идеи в большом потоке информации очень легко потеряться. Для упрощения генерации идеи лучше всего найти ограничения: желание разработчика и его аудитории писать в любимом жанре, выбор востребованного сеттинга, возможности </w:t></w:r><w:r><w:t xml:space="preserve">компьютерных </w:t></w:r><w:r><w:t>устройств и существующих технологий</w:t></w:r><w:r><w:t xml:space="preserve"> для разработки игр. </w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve">Часто идеи игр могут пересекаться или даже повторять друг друга, но это не делает игру хуже, поскольку гораздо важнее реализация. Возьму в качестве примера две игры в жанре королевской битвы: </w:t></w:r><w:r><w:t>Fortnite</w:t></w:r><w:r><w:t xml:space="preserve"> и </w:t></w:r><w:r><w:t>PUBG</w:t></w:r><w:r><w:t xml:space="preserve">. Обе игры нацелены на выживание среди 100 человек, находясь на одной карте, изначально имея стартовый набор. Отличается тем, что в </w:t></w:r><w:r><w:t>Fortnite</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">игрок делает упор на такую игровую механику, как строительство, создавая укрытия и </w:t></w:r><w:r><w:t>убегая от врагов, в</w:t></w:r><w:r><w:t xml:space="preserve"> то время как в </w:t></w:r><w:r><w:t>PUBG</w:t></w:r><w:r><w:t xml:space="preserve"> важнее сосредоточиться на стратегии </w:t></w:r><w:r><w:t>и скрытности.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Fortnite</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>выполнен в мультяшном стиле</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>PUBG</w:t></w:r><w:r><w:t xml:space="preserve"> в</w:t></w:r><w:r><w:t xml:space="preserve"> реалистичном. Обе игры имеют большие целевые аудитории и уникальны в своем жанре. </w:t></w:r></w:p><w:p><w:r><w:tab/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
